# Update Ingame character data: add two new characters (chr_deer, chr_mouse)
# to the "Character" sheet, as rows 8 and 9, mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Character")

# New row 8: id=4, name=chr_deer
$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(8, 2).Value = "chr_deer"
$ws.Cells.Item(8, 3).Value = 400
$ws.Cells.Item(8, 4).Value = 100
$ws.Cells.Item(8, 5).Value = 200
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2
$ws.Cells.Item(8, 8).Value = 3

# New row 9: id=5, name=chr_mouse
$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = "chr_mouse"
$ws.Cells.Item(9, 3).Value = 400
$ws.Cells.Item(9, 4).Value = 100
$ws.Cells.Item(9, 5).Value = 200
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2
$ws.Cells.Item(9, 8).Value = 3

# Copy styling from the row above (row 7) onto the new rows so formatting matches.
$ws.Range("A7:H7").Copy() | Out-Null
$ws.Range("A8:H9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Make the Character sheet the active sheet/tab, with B7 selected, matching the diff.
$ws.Activate()
$ws.Range("B7").Select() | Out-Null
